$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 118: ID 117, own (non-shared) DATE formula, B2B / Whatsapp / Nao Respondeu ---
$ws.Rows.Item(118).Insert()
$ws.Range("A118").Value = 117
$ws.Range("B118").Formula = "=DATE(2025,6,24)"
$ws.Range("C118").Value = "B2B"
$ws.Range("D118").Value = "Whatsapp"
$ws.Range("E118").Value = "Não Respondeu"

# --- Rows 119-140: ID 118-139, single shared DATE formula filled down together ---
$ws.Range("A119:A140").EntireRow.Insert()
$ws.Range("B119:B140").Formula = "=DATE(2025,6,24)"

# Rows 119-122: B2B / Whatsapp / Nao Respondeu
for ($i = 119; $i -le 122; $i++) {
    $ws.Cells.Item($i, 1).Value = $i - 1
    $ws.Cells.Item($i, 3).Value = "B2B"
    $ws.Cells.Item($i, 4).Value = "Whatsapp"
    $ws.Cells.Item($i, 5).Value = "Não Respondeu"
}

# Rows 123-140: B2B / Linkedin (no value in column E)
for ($i = 123; $i -le 140; $i++) {
    $ws.Cells.Item($i, 1).Value = $i - 1
    $ws.Cells.Item($i, 3).Value = "B2B"
    $ws.Cells.Item($i, 4).Value = "Linkedin"
}

# --- Row 141: blank trailing row, only the ID column keeps its number style ---
$ws.Rows.Item(141).Insert()
$ws.Range("B141").Clear()

# E2 no longer carries the stray explicit style that used to mark the old
# "last edited" cell - it now sits on E122 instead.
$ws.Range("E2").ClearFormats()

$ws.Range("E122").Select()
